$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RawDataPath for the ALUM_2020-IUCNGET row (row 9) to point to the
# new prerelease2 ABARES Land Use of Australia 2020-21 dataset.
$ws.Range("B9").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\inputs\raw\Land_use_of_Australia\ABARES_Land_use_of_Australia_2020_21_prerelease2_20240724\ABARES_Land_use_of_Australia_2020_21_prerelease2_20240724\NLUM_v7p2_ALUMV8_250m_2020_21_alb.tif"

# Move the active selection to B9, matching the saved workbook view state.
$ws.Range("B9").Select() | Out-Null
